$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; temporarily unprotect so the locked
# data cells can be updated, then restore protection afterwards.
$ws.Unprotect()

# Update the "as of" date in the confidentiality notice.
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-22 for illustrative purposes only and are subject to change."
# Re-fit the row height after the text update so it matches its original (auto) height.
$ws.Rows.Item(9).AutoFit()

# Refresh the Weight / Percent Change figures for each sector row.
$ws.Range("D2").Value = 0.2500946830959018
$ws.Range("E2").Value = -0.003589743589743399

$ws.Range("D3").Value = 0.2501715905651857
$ws.Range("E3").Value = -0.01724137931034475

$ws.Range("D4").Value = 0.2463866051778686
$ws.Range("E4").Value = 0.01751543209876538

$ws.Range("D5").Value = 0.2533471211610439
$ws.Range("E5").Value = 0.006630420368651579

$ws.Range("E6").Value = 0.0007842866944447469

$ws.Protect()
